$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.81798233333333
$ws.Range("H2").Value = 62.453947
$ws.Range("I2").Value = 0.8242653639952813
$ws.Range("J2").Value = 0.8242653639952813
$ws.Range("M2").Value = 58.95713633333333
$ws.Range("N2").Value = 176.871409
$ws.Range("O2").Value = 0.4863146960083892
$ws.Range("P2").Value = 0.4863146960083893
$ws.Range("Q2").Value = 1227.368622611258
$ws.Range("R2").Value = 11046.31760350132
$ws.Range("S2").Value = 0.4008523599216095
$ws.Range("T2").Value = 0.4008523599216096
$ws.Range("G3").Value = 20.81798233333333
$ws.Range("H3").Value = 62.453947
$ws.Range("I3").Value = 0.8242653639952813
$ws.Range("J3").Value = 0.8242653639952813
$ws.Range("O3").Value = 0.07416766570679004
$ws.Range("P3").Value = 0.07416766570679005
$ws.Range("Q3").Value = 187.1855126896367
$ws.Range("R3").Value = 1684.66961420673
$ws.Range("S3").Value = 0.06113383797048763
$ws.Range("T3").Value = 0.06113383797048765
$ws.Range("G4").Value = 20.81798233333333
$ws.Range("H4").Value = 62.453947
$ws.Range("I4").Value = 0.8242653639952813
$ws.Range("J4").Value = 0.8242653639952813
$ws.Range("M4").Value = 42.51661933333333
$ws.Range("N4").Value = 127.549858
$ws.Range("O4").Value = 0.3507032073181665
$ws.Range("P4").Value = 0.3507032073181665
$ws.Range("Q4").Value = 885.1102301543917
$ws.Range("R4").Value = 7965.992071389526
$ws.Range("S4").Value = 0.289072506834421
$ws.Range("T4").Value = 0.2890725068344211
$ws.Range("G5").Value = 20.81798233333333
$ws.Range("H5").Value = 62.453947
$ws.Range("I5").Value = 0.8242653639952813
$ws.Range("J5").Value = 0.8242653639952813
$ws.Range("M5").Value = 10.76719366666667
$ws.Range("N5").Value = 32.301581
$ws.Range("O5").Value = 0.0888144309666542
$ws.Range("P5").Value = 0.08881443096665421
$ws.Range("Q5").Value = 224.1512475322452
$ws.Range("R5").Value = 2017.361227790207
$ws.Range("S5").Value = 0.07320665926876301
$ws.Range("T5").Value = 0.07320665926876302
$ws.Range("G6").Value = 0.7925996666666667
$ws.Range("I6").Value = 0.03138212158540782
$ws.Range("J6").Value = 0.03138212158540782
$ws.Range("M6").Value = 58.95713633333333
$ws.Range("N6").Value = 176.871409
$ws.Range("O6").Value = 0.4863146960083892
$ws.Range("P6").Value = 0.4863146960083893
$ws.Range("Q6").Value = 46.72940660542123
$ws.Range("R6").Value = 420.564659448791
$ws.Range("S6").Value = 0.01526158691890591
$ws.Range("T6").Value = 0.01526158691890592
$ws.Range("G7").Value = 0.7925996666666667
$ws.Range("I7").Value = 0.03138212158540782
$ws.Range("J7").Value = 0.03138212158540782
$ws.Range("O7").Value = 0.07416766570679004
$ws.Range("P7").Value = 0.07416766570679005
$ws.Range("Q7").Value = 7.126683680823333
$ws.Range("R7").Value = 64.14015312741
$ws.Range("S7").Value = 0.002327538702916367
$ws.Range("T7").Value = 0.002327538702916367
$ws.Range("G8").Value = 0.7925996666666667
$ws.Range("I8").Value = 0.03138212158540782
$ws.Range("J8").Value = 0.03138212158540782
$ws.Range("M8").Value = 42.51661933333333
$ws.Range("N8").Value = 127.549858
$ws.Range("O8").Value = 0.3507032073181665
$ws.Range("P8").Value = 0.3507032073181665
$ws.Range("Q8").Value = 33.69865831139356
$ws.Range("R8").Value = 303.287924802542
$ws.Range("S8").Value = 0.01100581069245118
$ws.Range("T8").Value = 0.01100581069245119
$ws.Range("G9").Value = 0.7925996666666667
$ws.Range("I9").Value = 0.03138212158540782
$ws.Range("J9").Value = 0.03138212158540782
$ws.Range("M9").Value = 10.76719366666667
$ws.Range("N9").Value = 32.301581
$ws.Range("O9").Value = 0.0888144309666542
$ws.Range("P9").Value = 0.08881443096665421
$ws.Range("Q9").Value = 8.534074111135444
$ws.Range("R9").Value = 76.80666700021899
$ws.Range("S9").Value = 0.002787185271134351
$ws.Range("T9").Value = 0.002787185271134352
$ws.Range("G10").Value = 1.536855
$ws.Range("H10").Value = 4.610564999999999
$ws.Range("I10").Value = 0.06085010188305478
$ws.Range("J10").Value = 0.06085010188305479
$ws.Range("M10").Value = 58.95713633333333
$ws.Range("N10").Value = 176.871409
$ws.Range("O10").Value = 0.4863146960083892
$ws.Range("P10").Value = 0.4863146960083893
$ws.Range("Q10").Value = 90.60856975956499
$ws.Range("R10").Value = 815.4771278360848
$ws.Range("S10").Value = 0.0295922987993373
$ws.Range("T10").Value = 0.0295922987993373
$ws.Range("G11").Value = 1.536855
$ws.Range("H11").Value = 4.610564999999999
$ws.Range("I11").Value = 0.06085010188305478
$ws.Range("J11").Value = 0.06085010188305479
$ws.Range("O11").Value = 0.07416766570679004
$ws.Range("P11").Value = 0.07416766570679005
$ws.Range("Q11").Value = 13.81867783815
$ws.Range("R11").Value = 124.36810054335
$ws.Range("S11").Value = 0.004513110014686522
$ws.Range("T11").Value = 0.004513110014686523
$ws.Range("G12").Value = 1.536855
$ws.Range("H12").Value = 4.610564999999999
$ws.Range("I12").Value = 0.06085010188305478
$ws.Range("J12").Value = 0.06085010188305479
$ws.Range("M12").Value = 42.51661933333333
$ws.Range("N12").Value = 127.549858
$ws.Range("O12").Value = 0.3507032073181665
$ws.Range("P12").Value = 0.3507032073181665
$ws.Range("Q12").Value = 65.34187900552999
$ws.Range("R12").Value = 588.0769110497699
$ws.Range("S12").Value = 0.02134032589602451
$ws.Range("T12").Value = 0.02134032589602452
$ws.Range("G13").Value = 1.536855
$ws.Range("H13").Value = 4.610564999999999
$ws.Range("I13").Value = 0.06085010188305478
$ws.Range("J13").Value = 0.06085010188305479
$ws.Range("M13").Value = 10.76719366666667
$ws.Range("N13").Value = 32.301581
$ws.Range("O13").Value = 0.0888144309666542
$ws.Range("P13").Value = 0.08881443096665421
$ws.Range("Q13").Value = 16.547615422585
$ws.Range("R13").Value = 148.928538803265
$ws.Range("S13").Value = 0.005404367173006444
$ws.Range("T13").Value = 0.005404367173006445
$ws.Range("G14").Value = 2.108971
$ws.Range("H14").Value = 6.326912999999999
$ws.Range("I14").Value = 0.08350241253625613
$ws.Range("J14").Value = 0.08350241253625615
$ws.Range("M14").Value = 58.95713633333333
$ws.Range("N14").Value = 176.871409
$ws.Range("O14").Value = 0.4863146960083892
$ws.Range("P14").Value = 0.4863146960083893
$ws.Range("Q14").Value = 124.3388907700463
$ws.Range("R14").Value = 1119.050016930417
$ws.Range("S14").Value = 0.04060845036853651
$ws.Range("T14").Value = 0.04060845036853653
$ws.Range("G15").Value = 2.108971
$ws.Range("H15").Value = 6.326912999999999
$ws.Range("I15").Value = 0.08350241253625613
$ws.Range("J15").Value = 0.08350241253625615
$ws.Range("O15").Value = 0.07416766570679004
$ws.Range("P15").Value = 0.07416766570679005
$ws.Range("Q15").Value = 18.96287601563
$ws.Range("R15").Value = 170.66588414067
$ws.Range("S15").Value = 0.006193179018699519
$ws.Range("T15").Value = 0.006193179018699521
$ws.Range("G16").Value = 2.108971
$ws.Range("H16").Value = 6.326912999999999
$ws.Range("I16").Value = 0.08350241253625613
$ws.Range("J16").Value = 0.08350241253625615
$ws.Range("M16").Value = 42.51661933333333
$ws.Range("N16").Value = 127.549858
$ws.Range("O16").Value = 0.3507032073181665
$ws.Range("P16").Value = 0.3507032073181665
$ws.Range("Q16").Value = 89.66631719203933
$ws.Range("R16").Value = 806.9968547283539
$ws.Range("S16").Value = 0.0292845638952697
$ws.Range("T16").Value = 0.02928456389526971
$ws.Range("G17").Value = 2.108971
$ws.Range("H17").Value = 6.326912999999999
$ws.Range("I17").Value = 0.08350241253625613
$ws.Range("J17").Value = 0.08350241253625615
$ws.Range("M17").Value = 10.76719366666667
$ws.Range("N17").Value = 32.301581
$ws.Range("O17").Value = 0.0888144309666542
$ws.Range("P17").Value = 0.08881443096665421
$ws.Range("Q17").Value = 22.70769919438366
$ws.Range("R17").Value = 204.369292749453
$ws.Range("S17").Value = 0.007416219253750401
$ws.Range("T17").Value = 0.007416219253750403